$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = '@'
$c.Value = '37.074.19'
$c.Style = 'Normal'
$c = $ws.Range("E2")
$c.NumberFormat = '@'
$c.Value = '  +0.11%  '
$c.Style = 'Normal'
$c = $ws.Range("D3")
$c.NumberFormat = '@'
$c.Value = '2.049.29'
$c.Style = 'Normal'
$c = $ws.Range("E3")
$c.NumberFormat = '@'
$c.Value = '  -0.53%  '
$c.Style = 'Normal'
$c = $ws.Range("E4")
$c.NumberFormat = '@'
$c.Value = '  -0.12%  '
$c.Style = 'Normal'
$c = $ws.Range("D5")
$c.NumberFormat = '@'
$c.Value = '248.85'
$c.Style = 'Normal'
$c = $ws.Range("E5")
$c.NumberFormat = '@'
$c.Value = '  -0.24%  '
$c.Style = 'Normal'
$c = $ws.Range("D6")
$c.NumberFormat = '@'
$c.Value = '0.668'
$c.Style = 'Normal'
$c = $ws.Range("E6")
$c.NumberFormat = '@'
$c.Value = '  -0.55%  '
$c.Style = 'Normal'
$c = $ws.Range("D7")
$c.NumberFormat = '@'
$c.Value = '59.20'
$c.Style = 'Normal'
$c = $ws.Range("E7")
$c.NumberFormat = '@'
$c.Value = '  +8.59%  '
$c.Style = 'Normal'
$c = $ws.Range("E8")
$c.NumberFormat = '@'
$c.Value = '  -0.01%  '
$c.Style = 'Normal'
$c = $ws.Range("E9")
$c.NumberFormat = '@'
$c.Value = '  +0.97%  '
$c.Style = 'Normal'
$c = $ws.Range("D10")
$c.NumberFormat = '@'
$c.Value = '0.0799'
$c.Style = 'Normal'
$c = $ws.Range("E10")
$c.NumberFormat = '@'
$c.Value = '  +0.54%  '
$c.Style = 'Normal'
$c = $ws.Range("E11")
$c.NumberFormat = '@'
$c.Value = '  +2.04%  '
$c.Style = 'Normal'
$c = $ws.Range("D12")
$c.NumberFormat = '@'
$c.Value = '15.97'
$c.Style = 'Normal'
$c = $ws.Range("E12")
$c.NumberFormat = '@'
$c.Value = '  +6.32%  '
$c.Style = 'Normal'
$c = $ws.Range("D13")
$c.NumberFormat = '@'
$c.Value = '2.347.06'
$c.Style = 'Normal'
$c = $ws.Range("E13")
$c.NumberFormat = '@'
$c.Value = '  -0.60%  '
$c.Style = 'Normal'
$c = $ws.Range("D14")
$c.NumberFormat = '@'
$c.Value = '0.833'
$c.Style = 'Normal'
$c = $ws.Range("E14")
$c.NumberFormat = '@'
$c.Value = '  +1.90%  '
$c.Style = 'Normal'
$c = $ws.Range("D15")
$c.NumberFormat = '@'
$c.Value = '5.72'
$c.Style = 'Normal'
$c = $ws.Range("E15")
$c.NumberFormat = '@'
$c.Value = '  +7.31%  '
$c.Style = 'Normal'
$c = $ws.Range("D16")
$c.NumberFormat = '@'
$c.Value = '2.045.12'
$c.Style = 'Normal'
$c = $ws.Range("E16")
$c.NumberFormat = '@'
$c.Value = '  -0.78%  '
$c.Style = 'Normal'
$c = $ws.Range("D17")
$c.NumberFormat = '@'
$c.Value = '18.49'
$c.Style = 'Normal'
$c = $ws.Range("E17")
$c.NumberFormat = '@'
$c.Value = '  +29.84%  '
$c.Style = 'Normal'
$c = $ws.Range("D18")
$c.NumberFormat = '@'
$c.Value = '37.032.63'
$c.Style = 'Normal'
$c = $ws.Range("E18")
$c.NumberFormat = '@'
$c.Value = '  +0.04%  '
$c.Style = 'Normal'
$c = $ws.Range("E19")
$c.NumberFormat = '@'
$c.Value = '  +2.63%  '
$c.Style = 'Normal'
$c = $ws.Range("D20")
$c.NumberFormat = '@'
$c.Value = '0.0₃0902'
$c.Style = 'Normal'
$c = $ws.Range("E20")
$c.NumberFormat = '@'
$c.Value = '  -3.13%  '
$c.Style = 'Normal'
$c = $ws.Range("E21")
$c.NumberFormat = '@'
$c.Value = '  +0.60%  '
$c.Style = 'Normal'
$c = $ws.Range("D22")
$c.NumberFormat = '@'
$c.Value = '237.87'
$c.Style = 'Normal'
$c = $ws.Range("E22")
$c.NumberFormat = '@'
$c.Value = '  +0.28%  '
$c.Style = 'Normal'
$c = $ws.Range("E23")
$c.NumberFormat = '@'
$c.Value = '  +0.04%  '
$c.Style = 'Normal'
$c = $ws.Range("E24")
$c.NumberFormat = '@'
$c.Value = '  -0.38%  '
$c.Style = 'Normal'
$c = $ws.Range("E25")
$c.NumberFormat = '@'
$c.Value = '  +10.82%  '
$c.Style = 'Normal'
$c = $ws.Range("D26")
$c.NumberFormat = '@'
$c.Value = '9.48'
$c.Style = 'Normal'
$c = $ws.Range("E26")
$c.NumberFormat = '@'
$c.Value = '  +4.95%  '
$c.Style = 'Normal'
$c = $ws.Range("D27")
$c.NumberFormat = '@'
$c.Value = '168.81'
$c.Style = 'Normal'
$c = $ws.Range("E27")
$c.NumberFormat = '@'
$c.Value = '  -0.60%  '
$c.Style = 'Normal'
$c = $ws.Range("D28")
$c.NumberFormat = '@'
$c.Value = '20.09'
$c.Style = 'Normal'
$c = $ws.Range("E28")
$c.NumberFormat = '@'
$c.Value = '  -0.11%  '
$c.Style = 'Normal'
$c = $ws.Range("E29")
$c.NumberFormat = '@'
$c.Value = '  +0.31%  '
$c.Style = 'Normal'
$c = $ws.Range("E30")
$c.NumberFormat = '@'
$c.Value = '  +5.94%  '
$c.Style = 'Normal'
$c = $ws.Range("D31")
$c.NumberFormat = '@'
$c.Value = '4.79'
$c.Style = 'Normal'
$c = $ws.Range("E31")
$c.NumberFormat = '@'
$c.Value = '  +3.90%  '
$c.Style = 'Normal'
$c = $ws.Range("E32")
$c.NumberFormat = '@'
$c.Value = '  -0.50%  '
$c.Style = 'Normal'
$c = $ws.Range("E33")
$c.NumberFormat = '@'
$c.Value = '  +2.09%  '
$c.Style = 'Normal'
$c = $ws.Range("D34")
$c.NumberFormat = '@'
$c.Value = '0.0890'
$c.Style = 'Normal'
$c = $ws.Range("E34")
$c.NumberFormat = '@'
$c.Value = '  -0.71%  '
$c.Style = 'Normal'
$c = $ws.Range("E35")
$c.NumberFormat = '@'
$c.Value = '  +0.05%  '
$c.Style = 'Normal'
$c = $ws.Range("E36")
$c.NumberFormat = '@'
$c.Value = '  -2.75%  '
$c.Style = 'Normal'
$c = $ws.Range("E37")
$c.NumberFormat = '@'
$c.Value = '  -1.52%  '
$c.Style = 'Normal'
$c = $ws.Range("E38")
$c.NumberFormat = '@'
$c.Value = '  +4.44%  '
$c.Style = 'Normal'
$c = $ws.Range("D39")
$c.NumberFormat = '@'
$c.Value = '1.34'
$c.Style = 'Normal'
$c = $ws.Range("E39")
$c.NumberFormat = '@'
$c.Value = '  -1.03%  '
$c.Style = 'Normal'
$c = $ws.Range("D40")
$c.NumberFormat = '@'
$c.Value = '3.09'
$c.Style = 'Normal'
$c = $ws.Range("E40")
$c.NumberFormat = '@'
$c.Value = '  +10.44%  '
$c.Style = 'Normal'
$c = $ws.Range("D41")
$c.NumberFormat = '@'
$c.Value = '5.09'
$c.Style = 'Normal'
$c = $ws.Range("E41")
$c.NumberFormat = '@'
$c.Value = '  +23.20%  '
$c.Style = 'Normal'
$c = $ws.Range("D42")
$c.NumberFormat = '@'
$c.Value = '17.67'
$c.Style = 'Normal'
$c = $ws.Range("E42")
$c.NumberFormat = '@'
$c.Value = '  -0.01%  '
$c.Style = 'Normal'
$c = $ws.Range("E43")
$c.NumberFormat = '@'
$c.Value = '  -1.27%  '
$c.Style = 'Normal'
$c = $ws.Range("E44")
$c.NumberFormat = '@'
$c.Value = '  -0.79%  '
$c.Style = 'Normal'
$c = $ws.Range("D45")
$c.NumberFormat = '@'
$c.Value = '96.97'
$c.Style = 'Normal'
$c = $ws.Range("E45")
$c.NumberFormat = '@'
$c.Value = '  -0.01%  '
$c.Style = 'Normal'
$c = $ws.Range("D46")
$c.NumberFormat = '@'
$c.Value = '2.52'
$c.Style = 'Normal'
$c = $ws.Range("E46")
$c.NumberFormat = '@'
$c.Value = '  +5.17%  '
$c.Style = 'Normal'
$c = $ws.Range("D47")
$c.NumberFormat = '@'
$c.Value = '1.288.77'
$c.Style = 'Normal'
$c = $ws.Range("E47")
$c.NumberFormat = '@'
$c.Value = '  -0.57%  '
$c.Style = 'Normal'
$c = $ws.Range("D48")
$c.NumberFormat = '@'
$c.Value = '3.83'
$c.Style = 'Normal'
$c = $ws.Range("E48")
$c.NumberFormat = '@'
$c.Value = '  -7.80%  '
$c.Style = 'Normal'
$c = $ws.Range("D49")
$c.NumberFormat = '@'
$c.Value = '2.88'
$c.Style = 'Normal'
$c = $ws.Range("E49")
$c.NumberFormat = '@'
$c.Value = '  -1.23%  '
$c.Style = 'Normal'
$c = $ws.Range("D50")
$c.NumberFormat = '@'
$c.Value = '6.80'
$c.Style = 'Normal'
$c = $ws.Range("E50")
$c.NumberFormat = '@'
$c.Value = '  -1.13%  '
$c.Style = 'Normal'
$c = $ws.Range("D51")
$c.NumberFormat = '@'
$c.Value = '2.229.07'
$c.Style = 'Normal'
